$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F4").Value = "Asset"
$ws.Range("G4").Value = "asset-instance-1"

$ws.Range("H5").Select()
